$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats constant used for PasteSpecial
$xlPasteFormats = -4122

# --- Row 11 (2025-10-10) ---------------------------------------------------
# Apply the "Neutral" built-in cell style (already present in the sheet,
# e.g. on B2) to the cells that should turn orange/neutral.
$ws.Range("B2").Copy()
$ws.Range("B11").PasteSpecial($xlPasteFormats)
$ws.Range("F11").PasteSpecial($xlPasteFormats)
$ws.Range("J11").PasteSpecial($xlPasteFormats)

# Apply the "Good" built-in cell style (already present, e.g. on C2) to the
# cells that should turn green/good. G11 already has this style.
$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("E11").PasteSpecial($xlPasteFormats)
$ws.Range("H11").PasteSpecial($xlPasteFormats)
$ws.Range("I11").PasteSpecial($xlPasteFormats)

# Apply the "Bad" built-in cell style (already present, e.g. on D2) to the
# cells that should turn red/bad.
$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial($xlPasteFormats)
$ws.Range("K11").PasteSpecial($xlPasteFormats)
$ws.Range("L11").PasteSpecial($xlPasteFormats)
$ws.Range("M11").PasteSpecial($xlPasteFormats)

# Fill in the values, in the same order the author typed them so that the
# shared-string table ends up in the same order.
$ws.Range("G11").Value = "Entraînement I"
$ws.Range("B11").Value = "Sage"
$ws.Range("C11").Value = "Stratagème"
$ws.Range("E11").Value = "Distrait"
$ws.Range("F11").Value = "Rebelle"
$ws.Range("H11").Value = "Expertise artisanale"

# --- Row 31 (2025-10-30) ----------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("I31").PasteSpecial($xlPasteFormats)
$ws.Range("I31").Value = "Initiation magique"

# --- back to row 11 ----------------------------------------------------------
$ws.Range("I11").Value = "Attaque brutale, Lancer brutal"
$ws.Range("J11").Value = "Brute"

$excel.CutCopyMode = $false
